$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a literal text value into a cell (leading apostrophe prevents
# Excel from reinterpreting numeric-looking strings like "0.610" or "60.26" as
# numbers, which would silently drop trailing zeros / change the cell type).
# Resetting Style back to "Normal" afterwards keeps the cell style index identical
# to the original (no visual/formatting change), matching the source workbook.
function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.Value = "'" + $value
    $r.Style = "Normal"
}

# --- Price (column D) updates ---
Set-TextValue "D2" "37.823.56"
Set-TextValue "D3" "2.036.51"
Set-TextValue "D5" "228.14"
Set-TextValue "D6" "0.610"
Set-TextValue "D7" "60.26"
Set-TextValue "D11" "0.103"
Set-TextValue "D12" "2.340.77"
Set-TextValue "D13" "14.48"
Set-TextValue "D14" "21.12"
Set-TextValue "D15" "0.760"
Set-TextValue "D16" "5.18"
Set-TextValue "D17" "2.035.71"
Set-TextValue "D18" "37.800.90"
Set-TextValue "D19" "69.80"
Set-TextValue "D20" "5.90"
Set-TextValue "D22" "223.99"
Set-TextValue "D25" "2.26"
Set-TextValue "D27" "167.65"
Set-TextValue "D29" "18.87"
Set-TextValue "D32" "2.21"
Set-TextValue "D33" "4.40"
Set-TextValue "D34" "0.0606"
Set-TextValue "D35" "4.50"
Set-TextValue "D36" "6.33"
Set-TextValue "D40" "17.77"
Set-TextValue "D41" "1.536.65"
Set-TextValue "D42" "0.0217"
Set-TextValue "D43" "96.35"
Set-TextValue "D44" "2.80"
Set-TextValue "D45" "0.0913"
Set-TextValue "D46" "1.10"
Set-TextValue "D47" "4.01"
Set-TextValue "D49" "7.16"
Set-TextValue "D51" "2.229.50"

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("E4").Value = "  +0.34%  "
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  +1.18%  "
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  -2.12%  "
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("E16").Value = "  -2.27%  "
$ws.Range("E17").Value = "  -0.73%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("E20").Value = "  -5.92%  "
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +1.02%  "
$ws.Range("E27").Value = "  +0.83%  "
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("E29").Value = "  -1.02%  "
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("E32").Value = "  +7.20%  "
$ws.Range("E33").Value = "  -3.48%  "
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  +2.56%  "
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("E38").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +4.93%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("E44").Value = "  -2.52%  "
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -2.59%  "
$ws.Range("E47").Value = "  -3.73%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("E51").Value = "  -0.99%  "
